# Updates "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto symbol list, refreshing the scraped quotes. Values are written
# with a leading apostrophe so Excel keeps them as literal text (matching
# the original inlineStr cells) instead of re-interpreting them as
# numbers/percentages and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.60"
$ws.Range("E2").Value = "'-1.13%"
$ws.Range("D3").Value = "'35.73"
$ws.Range("E3").Value = "'-1.68%"
$ws.Range("D4").Value = "'5.039"
$ws.Range("E4").Value = "'-1.16%"
$ws.Range("D5").Value = "'0.07901"
$ws.Range("E5").Value = "'-3.14%"
$ws.Range("D6").Value = "'1.848"
$ws.Range("E6").Value = "'-4.49%"
$ws.Range("D7").Value = "'4.104"
$ws.Range("E7").Value = "'-2.27%"
$ws.Range("D8").Value = "'7.784"
$ws.Range("E8").Value = "'0.11%"
$ws.Range("D9").Value = "'0.9196"
$ws.Range("E9").Value = "'-1.15%"
$ws.Range("D10").Value = "'0.1339"
$ws.Range("E10").Value = "'-3.30%"
$ws.Range("D11").Value = "'0.1903"
$ws.Range("E11").Value = "'-1.12%"
$ws.Range("D12").Value = "'0.09126"
$ws.Range("E12").Value = "'-1.28%"
$ws.Range("D13").Value = "'0.03469"
$ws.Range("E13").Value = "'-3.80%"
$ws.Range("D14").Value = "'0.09823"
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("D15").Value = "'0.001399"
$ws.Range("E15").Value = "'-1.03%"
$ws.Range("D16").Value = "'0.006141"
$ws.Range("E16").Value = "'4.34%"
$ws.Range("D17").Value = "'3.726"
$ws.Range("E17").Value = "'3.67%"
$ws.Range("D18").Value = "'3.344"
$ws.Range("E18").Value = "'12.12%"
$ws.Range("E19").Value = "'0.04%"
$ws.Range("D20").Value = "'0.1310"
$ws.Range("E20").Value = "'-0.95%"
$ws.Range("D21").Value = "'5.158"
$ws.Range("E21").Value = "'5.57%"
$ws.Range("D22").Value = "'0.2191"
$ws.Range("E22").Value = "'-9.04%"
$ws.Range("D23").Value = "'0.04404"
$ws.Range("E23").Value = "'-2.47%"
$ws.Range("D24").Value = "'0.001233"
$ws.Range("E24").Value = "'1.48%"
$ws.Range("D25").Value = "'0.004620"
$ws.Range("E25").Value = "'-5.34%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'4.72%"
$ws.Range("D27").Value = "'0.0004443"
$ws.Range("E27").Value = "'0.03%"
$ws.Range("D39").Value = "'0.01937"
$ws.Range("E39").Value = "'-3.76%"
$ws.Range("D40").Value = "'0.05089"
$ws.Range("E40").Value = "'2.92%"
$ws.Range("D41").Value = "'0.007607"
$ws.Range("E41").Value = "'-1.07%"
$ws.Range("E42").Value = "'-8.73%"
$ws.Range("D43").Value = "'0.1342"
$ws.Range("E43").Value = "'-3.00%"
$ws.Range("D44").Value = "'0.002150"
$ws.Range("E44").Value = "'2.27%"
$ws.Range("D45").Value = "'0.01017"
$ws.Range("E45").Value = "'-5.32%"
$ws.Range("D46").Value = "'0.00006164"
$ws.Range("E46").Value = "'-4.39%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.09%"
$ws.Range("D48").Value = "'64.96"
$ws.Range("E48").Value = "'0.45%"
$ws.Range("E49").Value = "'39.34%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.09%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.09%"
